$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.506.97"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "1.907.31"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("E5").Value = "  +10.06%  "
$ws.Range("D6").Formula = "'246.57"
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Formula = "'40.82"
$ws.Range("E8").Value = "  -2.51%  "
$ws.Range("D9").Formula = "'0.356"
$ws.Range("E9").Value = "  +4.40%  "
$ws.Range("D10").Formula = "'52.54"
$ws.Range("E10").Value = "  +8.05%  "
$ws.Range("D11").Formula = "'0.0728"
$ws.Range("E11").Value = "  +3.40%  "
$ws.Range("D12").Formula = "'0.0987"
$ws.Range("E12").Value = "  -1.06%  "
$ws.Range("D13").Value = "2.183.74"
$ws.Range("E13").Value = "  +0.03%  "
$ws.Range("D14").Formula = "'12.62"
$ws.Range("E14").Value = "  +2.25%  "
$ws.Range("D15").Formula = "'0.715"
$ws.Range("E15").Value = "  +2.54%  "
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").Formula = "'4.90"
$ws.Range("E16").Value = "  +2.05%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "1.904.98"
$ws.Range("E17").Value = "  +0.42%  "
$ws.Range("D18").Value = "35.482.14"
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("D19").Formula = "'73.19"
$ws.Range("E19").Value = "  +1.81%  "
$ws.Range("D20").Value = "0.0₃0826"
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").Formula = "'242.68"
$ws.Range("E21").Value = "  -0.20%  "
$ws.Range("D22").Formula = "'13.07"
$ws.Range("E22").Value = "  +2.98%  "
$ws.Range("D23").Formula = "'5.05"
$ws.Range("E23").Value = "  +4.15%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").Formula = "'2.31"
$ws.Range("E25").Value = "  +0.39%  "
$ws.Range("D26").Formula = "'2.33"
$ws.Range("E26").Value = "  +6.52%  "
$ws.Range("D27").Formula = "'169.34"
$ws.Range("E27").Value = "  -1.63%  "
$ws.Range("D28").Formula = "'8.68"
$ws.Range("E28").Value = "  +1.44%  "
$ws.Range("D29").Formula = "'18.75"
$ws.Range("E29").Value = "  +4.37%  "
$ws.Range("E30").Value = "  +4.89%  "
$ws.Range("D32").Formula = "'4.24"
$ws.Range("E32").Value = "  +3.10%  "
$ws.Range("D33").Formula = "'0.0576"
$ws.Range("E33").Value = "  +0.94%  "
$ws.Range("D34").Formula = "'4.22"
$ws.Range("E34").Value = "  +1.25%  "
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").Formula = "'1.88"
$ws.Range("E35").Value = "  +6.15%  "
$ws.Range("B36").Value = "BinanceUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D36").Formula = "'1.01"
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("D37").Formula = "'0.917"
$ws.Range("E37").Value = "  -6.13%  "
$ws.Range("D38").Formula = "'1.47"
$ws.Range("E38").Value = "  +9.75%  "
$ws.Range("D39").Formula = "'2.02"
$ws.Range("E39").Value = "  -0.70%  "
$ws.Range("D40").Formula = "'96.70"
$ws.Range("E40").Value = "  +5.22%  "
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").Formula = "'16.90"
$ws.Range("E41").Value = "  +8.35%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Formula = "'1.12"
$ws.Range("E42").Value = "  +1.09%  "
$ws.Range("E43").Value = "  +2.38%  "
$ws.Range("E44").Value = "  +2.29%  "
$ws.Range("D45").Value = "1.354.61"
$ws.Range("E45").Value = "  +0.40%  "
$ws.Range("E46").Value = "  +1.36%  "
$ws.Range("B47").Value = "MultiversX"
$ws.Range("C47").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D47").Formula = "'46.27"
$ws.Range("E47").Value = "  -8.05%  "
$ws.Range("B48").Value = "HuobiToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D48").Formula = "'2.42"
$ws.Range("E48").Value = "  +0.19%  "
$ws.Range("E49").Value = "  +1.15%  "
$ws.Range("D50").Formula = "'12.33"
$ws.Range("E50").Value = "  -4.23%  "
$ws.Range("D51").Formula = "'6.53"
$ws.Range("E51").Value = "  -1.58%  "
